$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price reading was scraped for "Achicoria" at the Vega Central
# Mapocho de Santiago market. It belongs at the top of the data block (row 9,
# right after the header + the one row that stays fixed at row 8), so push
# the existing data rows 9:54 down one row to make room at row 9.
$ws.Rows.Item(9).Insert()

# Populate the newly inserted row 9 with the new reading.
$ws.Range("A9").Value = 9
$ws.Range("B9").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C9").Value = "Metropolitana"
$ws.Range("D9").Value = 45035
$ws.Range("E9").Value = 13
$ws.Range("F9").Value = 100112010
$ws.Range("G9").Value = "Achicoria"
$ws.Range("H9").Value = "Sin especificar"
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 70
$ws.Range("K9").Value = 7000
$ws.Range("L9").Value = 7000
$ws.Range("M9").Value = 7000
$ws.Range("N9").Value = "$/caja 16 unidades"
$ws.Range("O9").Value = "Provincia de Quillota"
$ws.Range("P9").Value = 438
$ws.Range("Q9").Value = 16
$ws.Range("R9").Value = "Hortaliza"
